$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I0 and IF in columns I and J of header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy header style (bold, bordered, centered) from H1 onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the I0 and IF data values for rows 2-66
$IVals = @(5,4,5,6,7,6,6,6,5,4,7,9,7,7,8,7,9,7,7,7,7,7,6,7,9,9,8,10,8,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,8,9,9,9,8,9,9,9,9,9,9,9,9,6,6,7,8,8,5)
$JVals = @(5,5,6,6,7,7,7,7,5,4,7,9,7,7,8,8,9,7,8,7,7,7,6,7,9,9,8,10,8,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,8,9,9,9,8,9,9,9,9,10,9,9,9,7,6,7,8,8,5)

for ($i = 0; $i -lt $IVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $IVals[$i]
    $ws.Cells.Item($row, 10).Value = $JVals[$i]
}
